$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format temporarily so numeric-looking strings
# like "1.00" or "13.40" are not auto-converted to numbers by Excel,
# preserving their original inline-string representation.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "51.588.56"
$ws.Range("E2").Value = "  +5.72%  "
$ws.Range("D3").Value = "2.739.63"
$ws.Range("E3").Value = "  +4.80%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "115.87"
$ws.Range("E5").Value = "  +5.20%  "
$ws.Range("D6").Value = "332.67"
$ws.Range("E6").Value = "  +3.29%  "
$ws.Range("E7").Value = "  +1.91%  "
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +5.06%  "
$ws.Range("D10").Value = "41.36"
$ws.Range("E10").Value = "  +4.95%  "
$ws.Range("E11").Value = "  +5.41%  "
$ws.Range("D12").Value = "20.08"
$ws.Range("E12").Value = "  +2.25%  "
$ws.Range("E13").Value = "  +2.75%  "
$ws.Range("D14").Value = "7.53"
$ws.Range("E14").Value = "  +4.50%  "
$ws.Range("D15").Value = "3.169.76"
$ws.Range("E15").Value = "  +4.80%  "
$ws.Range("D16").Value = "2.742.71"
$ws.Range("E16").Value = "  +4.67%  "
$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D17").Value = "0.873"
$ws.Range("E17").Value = "  +1.65%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "51.483.55"
$ws.Range("E18").Value = "  +5.56%  "
$ws.Range("E19").Value = "  +3.92%  "
$ws.Range("D20").Value = "13.40"
$ws.Range("D21").Value = "6.81"
$ws.Range("E21").Value = "  +2.09%  "
$ws.Range("E22").Value = "  +2.80%  "
$ws.Range("D23").Value = "278.49"
$ws.Range("E23").Value = "  +3.70%  "
$ws.Range("D24").Value = "69.04"
$ws.Range("E24").Value = "  +0.58%  "
$ws.Range("D25").Value = "2.64"
$ws.Range("E25").Value = "  +4.25%  "
$ws.Range("D26").Value = "26.60"
$ws.Range("E26").Value = "  +2.36%  "
$ws.Range("D28").Value = "10.15"
$ws.Range("E28").Value = "  +1.59%  "
$ws.Range("E29").Value = "  +0.19%  "
$ws.Range("D30").Value = "0.139"
$ws.Range("E30").Value = "  +2.18%  "
$ws.Range("D31").Value = "34.73"
$ws.Range("E31").Value = "  -0.16%  "
$ws.Range("D32").Value = "49.92"
$ws.Range("E32").Value = "  +1.12%  "
$ws.Range("D33").Value = "5.53"
$ws.Range("E33").Value = "  +0.90%  "
$ws.Range("E34").Value = "  +2.66%  "
$ws.Range("E35").Value = "  -0.12%  "
$ws.Range("D36").Value = "18.95"
$ws.Range("E36").Value = "  +0.75%  "
$ws.Range("E37").Value = "  -1.26%  "
$ws.Range("E38").Value = "  +1.90%  "
$ws.Range("D39").Value = "3.17"
$ws.Range("E39").Value = "  +2.25%  "
$ws.Range("D40").Value = "127.35"
$ws.Range("E40").Value = "  +0.56%  "
$ws.Range("D41").Value = "23.16"
$ws.Range("E41").Value = "  +4.97%  "
$ws.Range("D42").Value = "0.0345"
$ws.Range("E42").Value = "  +9.05%  "
$ws.Range("E43").Value = "  +7.71%  "
$ws.Range("E44").Value = "  +2.16%  "
$ws.Range("D45").Value = "2.39"
$ws.Range("E45").Value = "  +12.34%  "
$ws.Range("D46").Value = "2.088.08"
$ws.Range("E46").Value = "  +1.29%  "
$ws.Range("D47").Value = "3.31"
$ws.Range("E47").Value = "  +2.63%  "
$ws.Range("E48").Value = "  +3.47%  "
$ws.Range("E49").Value = "  +6.93%  "
$ws.Range("D50").Value = "8.90"
$ws.Range("E50").Value = "  +0.73%  "
$ws.Range("D51").Value = "59.57"
$ws.Range("E51").Value = "  +2.08%  "

# Restore default styling on column D (no explicit number format)
$ws.Range("D2:D51").Style = "Normal"

